$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("C32").Value = 5
$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 5

$ws.Range("F10").Select()
